$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.944.67'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '2.529.39'
$ws.Range("E3").Value = '  +3.45%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'537.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = "'143.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.55%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'0.571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").Value = '2.526.67'
$ws.Range("E9").Value = '  +2.81%  '
$ws.Range("D10").Value = "'0.0994"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").Value = "'5.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").Value = '2.941.36'
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").Value = "'23.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.61%  '
$ws.Range("D16").Value = '58.868.25'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '2.515.37'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").Value = "'11.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = "'4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.19%  '
$ws.Range("D21").Value = "'322.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").Value = "'5.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("D24").Value = "'61.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.89%  '
$ws.Range("D25").Value = "'0.436"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.58%  '
$ws.Range("E26").Value = '  +1.14%  '
$ws.Range("D27").Value = '2.614.24'
$ws.Range("E27").Value = '  +2.36%  '
$ws.Range("E28").Value = '  +1.89%  '
$ws.Range("D29").Value = "'7.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").Value = "'6.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.96%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0₃0765'
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("E33").Value = '  -8.38%  '
$ws.Range("D35").Value = "'158.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("E36").Value = '  +5.91%  '
$ws.Range("D37").Value = "'18.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("D38").Value = "'4.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("E39").Value = '  -6.69%  '
$ws.Range("D40").Value = "'36.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = "'5.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.39%  '
$ws.Range("D42").Value = "'296.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.19%  '
$ws.Range("D43").Value = "'3.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("D44").Value = "'0.808"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = "'0.603"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.96%  '
$ws.Range("D47").Value = "'10.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("D48").Value = "'124.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.61%  '
$ws.Range("D49").Value = "'0.0929"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("D50").Value = "'18.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("B51").Value = 'Hedera'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D51").Value = "'0.0512"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.80%  '
